# Column A on the "Card5" sheet stores the card number as text (rows 3-7 and
# 9-13 currently hold "2"). The commit changes those text values to "5"
# while leaving every other cell / style untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card5")

$targetCells = @("A3", "A4", "A5", "A6", "A7", "A9", "A10", "A11", "A12", "A13")
foreach ($addr in $targetCells) {
    $cell = $ws.Range($addr)
    # Force text (not numeric) storage, matching the rest of the column, then
    # restore the default "Normal" style so no formatting side effects remain.
    $cell.NumberFormat = "@"
    $cell.Value = "5"
    $cell.Style = "Normal"
}
